$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or "585.66"
# are not silently re-interpreted as numbers, matching the original inline-string data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "65.818.81"
$ws.Range("E2").Value = "  +6.79%  "

$ws.Range("D3").Value = "3.017.78"
$ws.Range("E3").Value = "  +4.40%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "585.66"
$ws.Range("E5").Value = "  +3.17%  "

$ws.Range("D6").Value = "155.63"
$ws.Range("E6").Value = "  +8.74%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.013.89"
$ws.Range("E8").Value = "  +4.27%  "

$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  +2.79%  "

$ws.Range("E11").Value = "  +6.17%  "

$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  +5.44%  "

$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +7.51%  "

$ws.Range("D14").Value = "34.46"
$ws.Range("E14").Value = "  +8.15%  "

$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "65.879.06"
$ws.Range("E16").Value = "  +6.93%  "

$ws.Range("D17").Value = "3.518.85"
$ws.Range("E17").Value = "  +4.43%  "

$ws.Range("D18").Value = "6.98"
$ws.Range("E18").Value = "  +6.61%  "

$ws.Range("D19").Value = "3.023.42"
$ws.Range("E19").Value = "  +4.53%  "

$ws.Range("D20").Value = "463.50"
$ws.Range("E20").Value = "  +7.11%  "

$ws.Range("D21").Value = "13.81"
$ws.Range("E21").Value = "  +5.78%  "

$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  +4.16%  "

$ws.Range("D23").Value = "7.38"
$ws.Range("E23").Value = "  +8.16%  "

$ws.Range("D24").Value = "82.18"
$ws.Range("E24").Value = "  +3.45%  "

$ws.Range("D25").Value = "12.58"
$ws.Range("E25").Value = "  +5.35%  "

$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  +12.43%  "

$ws.Range("D27").Value = "10.80"
$ws.Range("E27").Value = "  +9.06%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "2.43"
$ws.Range("E29").Value = "  +19.55%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "7.99"
$ws.Range("E30").Value = "  +14.43%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  +4.47%  "

$ws.Range("E33").Value = "  +5.93%  "

$ws.Range("D34").Value = "27.07"
$ws.Range("E34").Value = "  +6.07%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  +4.09%  "

$ws.Range("E37").Value = "  +7.88%  "

$ws.Range("D38").Value = "2.18"
$ws.Range("E38").Value = "  +13.06%  "

$ws.Range("E39").Value = "  +9.17%  "

$ws.Range("D40").Value = "49.35"
$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").Value = "45.19"
$ws.Range("E41").Value = "  +14.76%  "

$ws.Range("E42").Value = "  +8.01%  "

$ws.Range("D43").Value = "0.301"
$ws.Range("E43").Value = "  +13.37%  "

$ws.Range("D44").Value = "8.50"
$ws.Range("E44").Value = "  +3.65%  "

$ws.Range("D45").Value = "393.09"
$ws.Range("E45").Value = "  +13.36%  "

$ws.Range("D46").Value = "2.801.14"
$ws.Range("E46").Value = "  +4.24%  "

$ws.Range("D47").Value = "0.0355"
$ws.Range("E47").Value = "  +5.91%  "

$ws.Range("D48").Value = "134.83"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D50").Value = "23.85"
$ws.Range("E50").Value = "  +10.90%  "

$ws.Range("E51").Value = "  +4.23%  "

# Restore default (Normal) style on the price column so no stray number format remains
# attached to the cells (only the underlying value type/content should change).
$priceRange.Style = "Normal"
